# Unit test of libdoc
# Adds three new card entries to the "Misc" sheet:
#   - Yasiel Puig Heritage (row 51)
#   - Darren Daulton 96 UD Collectors Choice (row 336)
#   - Joel Embiid (new row 489, existing rows 489+ shift down by one)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# Row 336 (currently blank) - Darren Daulton 96 UD Collectors Choice
$ws.Range("B336").Value = "https://4.bp.blogspot.com/-GDHmyDkvlmk/WIAUdWK5pPI/AAAAAAAAzAI/m2b6nOHzGnE6vfLXKO8vTWCGomgmzcJEACLcB/s1600/PhotoScan%2B%25282%2529.jpg"
$ws.Range("A336").Value = "Darren Daulton 96 UD Collectors Choice"

# Row 51 (currently blank) - Yasiel Puig Heritage
$ws.Range("A51").Value = "Yasiel Puig Heritage"
$ws.Range("B51").Value = "https://4.bp.blogspot.com/-PlJjpeV-jeM/WIuu_dm_EwI/AAAAAAABr7s/29D9dZUfkNgmy6RHmvrn8fdBYqs0jX31gCLcB/s1600/julie-4.jpg"

# Insert a new row at 489 (shifts the NBA section rows down by one) - Joel Embiid
$ws.Rows.Item(489).Insert()
$ws.Range("B489").Value = "http://www.sportscollectorsdaily.com/wp-content/uploads/2017/01/Joel-Embiid-National-Treasures-relic.jpg"
$ws.Range("A489").Value = "Joel Embiid"

# Restore the view state recorded in the saved workbook
$ws.Application.ActiveWindow.ScrollRow = 415
$ws.Range("G490").Select()
